# Logged Week 15 and simulated Week 16
# Update Rushing and Receiving sheets with the latest per-player stats.

$wb = $excel.ActiveWorkbook

$rushing = $wb.Worksheets.Item("Rushing")
$receiving = $wb.Worksheets.Item("Receiving")

# --- Rushing sheet updates ---
# Row 2: T.Lawrence
$rushing.Range("C2").Value = 19
$rushing.Range("D2").Value = 11
$rushing.Range("E2").Value = 22
$rushing.Range("F2").Value = 10

# Row 4: J.Robinson
$rushing.Range("C4").Value = 90
$rushing.Range("D4").Value = 60
$rushing.Range("E4").Value = 10
$rushing.Range("F4").Value = 23

# Row 6: D.Ogunbowale
$rushing.Range("E6").Value = 1

# Row 7: L.Shenault
$rushing.Range("D7").Value = 7

# --- Receiving sheet updates ---
# Row 2: J.Robinson
$receiving.Range("C2").Value = 45
$receiving.Range("D2").Value = 30

# Row 4: D.Ogunbowale
$receiving.Range("C4").Value = 9
$receiving.Range("D4").Value = 6

# Row 5: M.Jones
$receiving.Range("C5").Value = 64
$receiving.Range("D5").Value = 46
$receiving.Range("E5").Value = 28
$receiving.Range("G5").Value = 10
$receiving.Range("H5").Value = 3

# Row 6: L.Shenault
$receiving.Range("C6").Value = 76
$receiving.Range("D6").Value = 51

# Row 9: T.Austin
$receiving.Range("C9").Value = 23
$receiving.Range("D9").Value = 14
$receiving.Range("E9").Value = 4
$receiving.Range("G9").Value = 3

# Row 10: L.Treadwell
$receiving.Range("C10").Value = 22
$receiving.Range("D10").Value = 16
$receiving.Range("E10").Value = 10
$receiving.Range("F10").Value = 4

# Row 12: C.Manhertz
$receiving.Range("C12").Value = 7
$receiving.Range("D12").Value = 5

# Row 13: J.O'Shaughnessy
$receiving.Range("C13").Value = 24
$receiving.Range("D13").Value = 16
$receiving.Range("E13").Value = 7
$receiving.Range("F13").Value = 4
